$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 49/50: EOS and Quant swapped position
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"

# D and E columns: numeric-looking / formatted text values must be forced as text
# (use a literal-string formula, then Copy + PasteSpecial values to "bake" it in
#  as text without altering cell styling or number formats).
$ws.Range("D2").Formula = '="28.420.54"'
$ws.Range("E2").Formula = '="  +0.73%  "'
$ws.Range("D3").Formula = '="1.872.77"'
$ws.Range("E3").Formula = '="  -0.47%  "'
$ws.Range("E4").Formula = '="  +0.89%  "'
$ws.Range("D5").Formula = '="316.67"'
$ws.Range("E6").Formula = '="  +1.05%  "'
$ws.Range("E7").Formula = '="  -0.47%  "'
$ws.Range("D8").Formula = '="0.3960"'
$ws.Range("E8").Formula = '="  +1.61%  "'
$ws.Range("D9").Formula = '="0.08436"'
$ws.Range("E9").Formula = '="  +1.09%  "'
$ws.Range("E10").Formula = '="  -0.94%  "'
$ws.Range("E11").Formula = '="  +0.66%  "'
$ws.Range("D12").Formula = '="6.248"'
$ws.Range("E12").Formula = '="  +0.51%  "'
$ws.Range("D13").Formula = '="1.888.12"'
$ws.Range("E13").Formula = '="  +0.85%  "'
$ws.Range("D14").Formula = '="20.45"'
$ws.Range("E14").Formula = '="  -0.43%  "'
$ws.Range("D15").Formula = '="7.229"'
$ws.Range("E15").Formula = '="  -0.24%  "'
$ws.Range("E16").Formula = '="  +0.91%  "'
$ws.Range("E17").Formula = '="  +0.81%  "'
$ws.Range("D18").Formula = '="90.95"'
$ws.Range("E18").Formula = '="  +0.00%  "'
$ws.Range("D19").Formula = '="0.06767"'
$ws.Range("E19").Formula = '="  +1.33%  "'
$ws.Range("D20").Formula = '="17.71"'
$ws.Range("E20").Formula = '="  -0.44%  "'
$ws.Range("E21").Formula = '="  +1.01%  "'
$ws.Range("D22").Formula = '="5.935"'
$ws.Range("E22").Formula = '="  -1.47%  "'
$ws.Range("D23").Formula = '="28.481.13"'
$ws.Range("E23").Formula = '="  +0.80%  "'
$ws.Range("D24").Formula = '="11.16"'
$ws.Range("E24").Formula = '="  +0.46%  "'
$ws.Range("D25").Formula = '="2.290"'
$ws.Range("E25").Formula = '="  +0.80%  "'
$ws.Range("D26").Formula = '="2.099.18"'
$ws.Range("E26").Formula = '="  +0.46%  "'
$ws.Range("D27").Formula = '="161.63"'
$ws.Range("E27").Formula = '="  +0.71%  "'
$ws.Range("D28").Formula = '="20.63"'
$ws.Range("E28").Formula = '="  -0.11%  "'
$ws.Range("D29").Formula = '="2.343"'
$ws.Range("E29").Formula = '="  -4.36%  "'
$ws.Range("D30").Formula = '="127.03"'
$ws.Range("E30").Formula = '="  +1.24%  "'
$ws.Range("D31").Formula = '="0.1053"'
$ws.Range("E31").Formula = '="  -0.57%  "'
$ws.Range("D32").Formula = '="1.036"'
$ws.Range("E32").Formula = '="  -0.15%  "'
$ws.Range("D33").Formula = '="5.756"'
$ws.Range("E33").Formula = '="  -1.81%  "'
$ws.Range("D34").Formula = '="3.641"'
$ws.Range("E34").Formula = '="  +0.83%  "'
$ws.Range("D35").Formula = '="0.02431"'
$ws.Range("E35").Formula = '="  -0.27%  "'
$ws.Range("D36").Formula = '="0.06460"'
$ws.Range("E36").Formula = '="  -1.40%  "'
$ws.Range("D37").Formula = '="0.2174"'
$ws.Range("E37").Formula = '="  -1.47%  "'
$ws.Range("D38").Formula = '="8.793"'
$ws.Range("E38").Formula = '="  -6.86%  "'
$ws.Range("E39").Formula = '="  +1.51%  "'
$ws.Range("D40").Formula = '="1.181"'
$ws.Range("E40").Formula = '="  -1.52%  "'
$ws.Range("D41").Formula = '="0.6372"'
$ws.Range("E41").Formula = '="  -1.72%  "'
$ws.Range("D42").Formula = '="4.981"'
$ws.Range("E42").Formula = '="  -0.42%  "'
$ws.Range("D43").Formula = '="11.20"'
$ws.Range("E43").Formula = '="  -0.10%  "'
$ws.Range("D44").Formula = '="0.6032"'
$ws.Range("E44").Formula = '="  -0.83%  "'
$ws.Range("D45").Formula = '="13.02"'
$ws.Range("E45").Formula = '="  -1.07%  "'
$ws.Range("D46").Formula = '="3.711"'
$ws.Range("E46").Formula = '="  +0.49%  "'
$ws.Range("D47").Formula = '="1.988"'
$ws.Range("E47").Formula = '="  -1.25%  "'
$ws.Range("D48").Formula = '="1.206"'
$ws.Range("E48").Formula = '="  -5.86%  "'
$ws.Range("D49").Formula = '="121.96"'
$ws.Range("E49").Formula = '="  +0.76%  "'
$ws.Range("D50").Formula = '="1.203"'
$ws.Range("E50").Formula = '="  -2.74%  "'
$ws.Range("D51").Formula = '="0.06840"'
$ws.Range("E51").Formula = '="  -0.97%  "'

# Bake the formulas into static text values
$targets = @("D2", "E2", "D3", "E3", "E4", "D5", "E6", "E7", "D8", "E8", "D9", "E9", "E10", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "E16", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($ref in $targets) {
    $c = $ws.Range($ref)
    $c.Copy()
    $c.PasteSpecial(-4163)
}
